$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("D9").NumberFormat = "0%"
$ws.Range("D9").Value = "75% (Làm được 3/4 chức năng)"
$ws.Range("E9").Value = "100% (10/06/2010)"
$ws.Range("E10").Select()
